# Refresh the cryptos price list (rows 2-51) with the latest scraped
# values: updated prices (col D), updated 1h volume deltas (col E), and
# the Hedera/OKB rows (35/36) that swapped rank positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat "@" (Text) is applied before writing any price that looks
# like a plain decimal (e.g. "573.16") so Excel stores it as a string
# instead of silently converting it to a number; Style is reset back to
# "Normal" right after so no stray style index is left on the cell.

$ws.Range('D2').Value = '67.787.71'
$ws.Range('E2').Value = '  -6.49%  '

$ws.Range('D3').Value = '3.687.99'
$ws.Range('E3').Value = '  -5.98%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.89%  '

$ws.Range('D7').Value = '3.682.34'
$ws.Range('E7').Value = '  -5.94%  '

$ws.Range('E8').Value = '  -8.41%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.997'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('E10').Value = '  -9.45%  '

$ws.Range('E11').Value = '  -13.43%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.89%  '

$ws.Range('E13').Value = '  -12.43%  '

$ws.Range('E14').Value = '  -9.42%  '

$ws.Range('D15').Value = '4.283.14'
$ws.Range('E15').Value = '  -5.79%  '

$ws.Range('D16').Value = '3.686.31'
$ws.Range('E16').Value = '  -6.10%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -9.09%  '

$ws.Range('E18').Value = '  -3.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.66%  '

$ws.Range('E20').Value = '  -9.01%  '

$ws.Range('D21').Value = '67.512.09'
$ws.Range('E21').Value = '  -6.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '403.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.71%  '

$ws.Range('E25').Value = '  -8.61%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.77%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.82%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.78'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.07%  '

$ws.Range('E30').Value = '  -8.19%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.35'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.12%  '

$ws.Range('E32').Value = '  -6.51%  '

$ws.Range('E33').Value = '  -10.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '611.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.42%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.68'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.57%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.115'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -15.92%  '

$ws.Range('D38').Value = '0.0₃0880'
$ws.Range('E38').Value = '  -10.29%  '

$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.392'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E41').Value = '  -0.23%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.135'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.67%  '

$ws.Range('E43').Value = '  +7.36%  '

$ws.Range('E44').Value = '  -10.63%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0430'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.31%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.22%  '

$ws.Range('D48').Value = '2.791.91'
$ws.Range('E48').Value = '  -1.23%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.132'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.54%  '
